$d = $word.ActiveDocument
$results = @()

# ---------------------------------------------------------------------------
# 1. Top-of-document "Date:" paragraph.
#    The paragraph is three runs: "Date:" / " " / "2017-03-01". Only the
#    third run's text should change to "2017-04-24". A plain
#    Find.Execute(..., Replace:=wdReplaceOne) (or setting Range.Text) ends up
#    deleting+re-typing through the space run too and merges it with the
#    date run, so instead we locate just the date text, delete it and
#    re-insert the new date -- that keeps the run immediately after the
#    space run separate, matching the original structure.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Text = "2017-03-01"
$rng.Find.Forward = $true
$rng.Find.MatchCase = $true
$rng.Find.MatchWildcards = $false
$found = $rng.Find.Execute()
if ($found) {
    $rng.Delete()
    $rng.InsertAfter("2017-04-24")
}
$results += "date-field: $found"

# ---------------------------------------------------------------------------
# 2. "Report rendered by ..." line in the session-info appendix -- single run,
#    safe to replace with a plain Find/Replace.
# ---------------------------------------------------------------------------
$r2 = $d.Content.Find.Execute(
    "Report rendered by koval_000 at 2017-03-01, 09:04 -0500", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "Report rendered by koval_000 at 2017-04-24, 11:36 -0400", 2)
$results += "rendered-by-line: $r2"

# ---------------------------------------------------------------------------
# 3. `sessionInfo()` "other attached packages" / "loaded via a namespace"
#    listing -- each of the five lines below is its own run, so they can be
#    replaced wholesale.
# ---------------------------------------------------------------------------
$r3 = $d.Content.Find.Execute(
    "[1] knitr_1.15.1    forestplot_1.7  checkmate_1.8.2 ggplot2_2.2.1   magrittr_1.5   ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[1] dplyr_0.5.0     forestplot_1.7  checkmate_1.8.2 ggplot2_2.2.1   magrittr_1.5    knitr_1.15.1   ",
    2)
$results += "pkg-line-1: $r3"

$r4 = $d.Content.Find.Execute(
    " [7] stringr_1.1.0    plyr_1.8.4       dplyr_0.5.0      tools_3.3.2      DT_0.2           gtable_0.2.0    ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " [7] stringr_1.1.0    plyr_1.8.4       tools_3.3.2      DT_0.2           gtable_0.2.0     DBI_0.5-1       ",
    2)
$results += "pkg-line-2: $r4"

$r5 = $d.Content.Find.Execute(
    "[13] plotrix_3.6-4    DBI_0.5-1        htmltools_0.3.5  yaml_2.1.14      lazyeval_0.2.0   assertthat_0.1  ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[13] htmltools_0.3.5  yaml_2.1.14      lazyeval_0.2.0   assertthat_0.1   rprojroot_1.2    digest_0.6.12   ",
    2)
$results += "pkg-line-3: $r5"

$r6 = $d.Content.Find.Execute(
    "[19] digest_0.6.12    rprojroot_1.2    tibble_1.2       readr_1.0.0      tidyr_0.6.1      htmlwidgets_0.8 ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[19] tibble_1.2       readr_1.0.0      tidyr_0.6.1      htmlwidgets_0.8  evaluate_0.10    haven_1.0.0     ",
    2)
$results += "pkg-line-4: $r6"

$r7 = $d.Content.Find.Execute(
    "[25] evaluate_0.10    rmarkdown_1.3    stringi_1.1.2    scales_0.4.1     backports_1.0.5  jsonlite_1.2    ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[25] rmarkdown_1.3    stringi_1.1.2    scales_0.4.1     backports_1.0.5 ",
    2)
$results += "pkg-line-5: $r7"

# ---------------------------------------------------------------------------
# Note: the diff also rewrites the <w:nsid> GUID-like values of two
# abstractNum list definitions in numbering.xml (990: 60f5057a -> 23594062;
# 991: b1b22793 -> dab7d30d). These identifiers aren't tied to any paragraph
# content (numId 1000/abstractNum 990 isn't even referenced by the body) and
# there's no Word object-model property that exposes or edits them -- they
# are regenerated internally by the external rendering pipeline (knitr /
# pandoc) each time the report is re-knit, not something Word automation can
# set. No reachable COM call corresponds to that hunk.
# ---------------------------------------------------------------------------

$results
